$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing question-cell formatting (bold, bordered, centered style)
# from A7 onto the new question rows A20:A26 before writing values.
$ws.Range("A7").Copy()
$ws.Range("A20:A26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @("WHICH PART IS THIS?", "filler"),
    @("ON WHICH PORT WILL THE WIRE OF TORQUE MOTOR GOES?", "M4"),
    @("how many castor wheels are inserted in one-arm soccer bot", "1"),
    @("WHICH PART IS THIS FROM MODEL?", "l-channel"),
    @("WHICH SIZE OF AXLE WE HAVE TO INSERT IN TORQUE MOTOR?", "5.5mm axle"),
    @("Can we use castor wheel on the back of the model (true/false)?", "True"),
    @("ON WHICH PART WE WILL INSERT THE ONE ARM ?", "square plate")
)

$row = 20
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]

    $answer = $pair[1]
    $cell = $ws.Cells.Item($row, 2)
    # Prefix with an apostrophe so values that look numeric/boolean
    # ("1", "True") are stored as literal text, matching the source data.
    $cell.Value = "'" + $answer
    # Drop the quote-prefix formatting Excel applies for text-forced entry,
    # restoring the cell to the default (unstyled) look used elsewhere in
    # column B.
    $cell.Style = "Normal"

    $row++
}
